$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '297.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.74%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.90'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.85%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.40%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07518'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.68%'
$ws.Range("B6").Value = 'GateToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.378'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.94%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.585'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.30%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9259'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.20%'
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.401'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.36%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1194'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.41%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1825'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.01%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08883'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.27%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.04063'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-6.18%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.1050'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.57%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001278'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.14%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005841'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.14%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.356'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.55%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.74%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.075'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '1.19%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1351'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-3.51%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '13.07%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.50%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.001266'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.35%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.003911'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.54%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02415'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '6.46%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05212'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.67%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.006305'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '6.77%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007796'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.46%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1327'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.28%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007407'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.86%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007271'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-12.22%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.2970'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.87%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006579'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '4.24%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.07%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.03170'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '23.87%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.004203'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002102'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.07%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002002'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.07%'
